# Adds a header row ("Input Sheet" / "Value") to every worksheet in the
# workbook. Cells are bold, centered horizontally, top-aligned vertically,
# and boxed with a thin border.
#
# The header format (bold + thin box border + center/top alignment) is
# built up ONCE on a scratch cell and then copied across to every target
# cell via PasteSpecial (format-only), one single cell at a time. That
# keeps the whole property combination in a single shared style record
# instead of spawning a fresh style per property assignment / per sheet.

$wb = $excel.ActiveWorkbook

$templateWs = $wb.Worksheets.Item(1)
$template = $templateWs.Range("ZZ1")
$template.Value = "x"
$template.Font.Bold = $true
$template.Borders.LineStyle = 1
$template.HorizontalAlignment = -4108   # xlCenter
$template.VerticalAlignment = -4160     # xlTop
$template.Copy()

foreach ($ws in $wb.Worksheets) {
    $ws.Range("A1").Value = "Input Sheet"
    $ws.Range("B1").Value = "Value"
    $ws.Range("A1").PasteSpecial(-4122)  # xlPasteFormats
    $ws.Range("B1").PasteSpecial(-4122)  # xlPasteFormats
}

$template.Clear()
